# "Generate Report for Handback" -- refresh the handoff/handback timestamp
# columns on each sheet of the report, as a new localization report run
# would produce.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the first (6cabd8c3...)
# row moves from 09:07:14 to 09:08:19.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-03 09:08:19"

# zh-cn sheet: handoff/handback datetimes for the first (6cabd8c3...) row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-03 09:08:14"
$wsZhCn.Range("K2").Value = "2016-09-03 09:08:31"

# de-de sheet: handoff datetime (shared with the Overview value) and the
# handback datetime for the first (6cabd8c3...) row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-03 09:08:19"
$wsDeDe.Range("K2").Value = "2016-09-03 09:08:38"
